$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers but must remain stored as
# text (matching the original inline-string "Price" column formatting).
# We briefly flip the number format to Text, assign the value, then restore
# the cell style so no lingering numeric formatting is left behind.
$textValues = @{
    'D5' = '253.31'
    'D6' = '0.635'
    'D7' = '72.03'
    'D9' = '0.643'
    'D10' = '41.14'
    'D11' = '59.59'
    'D12' = '0.0964'
    'D13' = '7.37'
    'D16' = '14.81'
    'D17' = '0.882'
    'D21' = '6.25'
    'D22' = '72.95'
    'D23' = '235.72'
    'D24' = '2.11'
    'D26' = '11.74'
    'D31' = '167.79'
    'D32' = '20.98'
    'D33' = '0.131'
    'D35' = '0.0787'
    'D36' = '0.124'
    'D37' = '28.76'
    'D38' = '4.70'
    'D42' = '5.98'
    'D43' = '12.44'
    'D44' = '64.66'
    'D45' = '4.95'
    'D47' = '8.94'
}

foreach ($ref in $textValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$ref]
    $cell.Style = "Normal"
}

# Remaining cell updates (text/link/percentage values that Excel will not
# reinterpret as numbers).
$updates = @{
    'D2' = '42.726.44'
    'E2' = '  +3.38%  '
    'D3' = '2.253.58'
    'E3' = '  +3.25%  '
    'E4' = '  +0.01%  '
    'E5' = '  -0.63%  '
    'E6' = '  +1.39%  '
    'E7' = '  +5.80%  '
    'E8' = '  -0.16%  '
    'E9' = '  +11.79%  '
    'E10' = '  +9.12%  '
    'E11' = '  +0.93%  '
    'E12' = '  +3.22%  '
    'E13' = '  +3.39%  '
    'E14' = '  +0.55%  '
    'D15' = '2.590.19'
    'E15' = '  +3.42%  '
    'B16' = 'Chainlink'
    'C16' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E16' = '  +2.47%  '
    'B17' = 'Polygon'
    'C17' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E17' = '  +1.17%  '
    'D18' = '2.254.71'
    'E18' = '  +0.49%  '
    'D19' = '42.735.87'
    'E19' = '  +3.66%  '
    'D20' = '0.0₃0980'
    'E20' = '  +2.85%  '
    'E21' = '  +1.34%  '
    'E22' = '  +1.66%  '
    'E23' = '  +1.73%  '
    'E24' = '  +4.25%  '
    'E25' = '  +0.87%  '
    'E26' = '  +0.08%  '
    'E27' = '  +0.21%  '
    'E28' = '  -3.30%  '
    'E29' = '  -1.76%  '
    'E30' = '  +2.00%  '
    'E31' = '  -0.38%  '
    'E32' = '  +1.84%  '
    'E33' = '  +12.60%  '
    'E34' = '  +12.32%  '
    'E35' = '  +5.10%  '
    'E36' = '  +0.74%  '
    'E37' = '  +9.15%  '
    'E38' = '  +1.71%  '
    'E39' = '  +0.55%  '
    'E40' = '  +6.42%  '
    'E41' = '  +4.47%  '
    'E42' = '  +5.45%  '
    'E43' = '  -0.46%  '
    'E44' = '  +0.81%  '
    'E45' = '  -2.26%  '
    'E46' = '  +0.73%  '
    'E47' = '  +3.46%  '
    'E48' = '  +1.34%  '
    'E49' = '  +4.96%  '
    'E50' = '  -0.48%  '
    'E51' = '  +3.73%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
